$wb = $excel.ActiveWorkbook

# --- Update the text note on "Hoja1" (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = [string][char]0x2705 + " 1000 Bs = 11.4 = 46470.52 pesos"
$newLine1 = [string][char]0x2705 + " 1000 Bs = 11.33 = 46359.12 pesos"

$oldLine2 = [string][char]0x2705 + " 46470.52 pesos = 11.34 = 965.24 Bs"
$newLine2 = [string][char]0x2705 + " 46359.12 pesos = 11.3 = 974.8 Bs"

$text = $ws1.Range("A1").Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$ws1.Range("A1").Value = $text

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 88.289
$ws2.Range("O10").Value = 4093

$ws2.Range("N12").Value = 4104.2
$ws2.Range("O12").Value = 86.3
